$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.641.67'
$ws.Range('E2').Value = '  -2.34%  '
$ws.Range('D3').Value = '2.975.49'
$ws.Range('E3').Value = '  +0.89%  '
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').Value = '558.91'
$ws.Range('E5').Value = '  -2.19%  '
$ws.Range('D6').Value = '132.44'
$ws.Range('E6').Value = '  +7.70%  '
$ws.Range('E7').Value = '  -0.14%  '
$ws.Range('D8').Value = '0.517'
$ws.Range('E8').Value = '  +4.28%  '
$ws.Range('D9').Value = '2.965.74'
$ws.Range('E9').Value = '  +1.03%  '
$ws.Range('D10').Value = '0.129'
$ws.Range('E10').Value = '  -0.94%  '
$ws.Range('D11').Value = '4.86'
$ws.Range('E11').Value = '  -4.01%  '
$ws.Range('D12').Value = '0.451'
$ws.Range('E12').Value = '  +4.22%  '
$ws.Range('D13').Value = '0.0000225'
$ws.Range('E13').Value = '  +2.29%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '33.10'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +2.76%  '
$ws.Range('E15').Value = '  +2.28%  '
$ws.Range('D16').Value = '3.467.74'
$ws.Range('E16').Value = '  +0.45%  '
$ws.Range('D17').Value = '6.85'
$ws.Range('E17').Value = '  +12.14%  '
$ws.Range('D18').Value = '2.980.58'
$ws.Range('E18').Value = '  +0.40%  '
$ws.Range('D19').Value = '58.662.93'
$ws.Range('E19').Value = '  -2.53%  '
$ws.Range('D20').Value = '423.09'
$ws.Range('E20').Value = '  -0.66%  '
$ws.Range('D21').Value = '13.24'
$ws.Range('E21').Value = '  +2.60%  '
$ws.Range('D22').Value = '0.687'
$ws.Range('E22').Value = '  +4.97%  '
$ws.Range('D23').Value = '7.03'
$ws.Range('E23').Value = '  +1.17%  '
$ws.Range('D24').Value = '13.07'
$ws.Range('D25').Value = '79.89'
$ws.Range('E25').Value = '  +2.24%  '
$ws.Range('E26').Value = '  +0.20%  '
$ws.Range('D27').Value = '0.999'
$ws.Range('E27').Value = '  -0.25%  '
$ws.Range('D28').Value = '2.51'
$ws.Range('E28').Value = '  -0.18%  '
$ws.Range('D29').Value = '7.64'
$ws.Range('E29').Value = '  +8.21%  '
$ws.Range('E30').Value = '  +8.56%  '
$ws.Range('D31').Value = '0.106'
$ws.Range('E31').Value = '  +16.15%  '
$ws.Range('D32').Value = '6.22'
$ws.Range('E32').Value = '  +3.07%  '
$ws.Range('D33').Value = '25.28'
$ws.Range('E33').Value = '  +1.09%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.70'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +3.65%  '
$ws.Range('B35').Value = 'Stacks'
$ws.Range('C35').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D35').Value = '2.14'
$ws.Range('E35').Value = '  -2.02%  '
$ws.Range('B36').Value = 'Mantle'
$ws.Range('C36').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D36').Value = '0.952'
$ws.Range('E36').Value = '  +1.84%  '
$ws.Range('B37').Value = 'PEPE'
$ws.Range('C37').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D37').Value = '0.0₃0697'
$ws.Range('E37').Value = '  +8.25%  '
$ws.Range('B38').Value = 'OKB'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '48.70'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -0.88%  '
$ws.Range('D39').Value = '8.47'
$ws.Range('E39').Value = '  +8.62%  '
$ws.Range('E40').Value = '  +10.54%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.110'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +2.02%  '
$ws.Range('D42').Value = '0.0353'
$ws.Range('E42').Value = '  -0.45%  '
$ws.Range('D43').Value = '381.34'
$ws.Range('E43').Value = '  +2.42%  '
$ws.Range('D44').Value = '2.665.70'
$ws.Range('E44').Value = '  +1.50%  '
$ws.Range('E45').Value = '  -0.01%  '
$ws.Range('D46').Value = '0.242'
$ws.Range('E46').Value = '  +4.25%  '
$ws.Range('D47').Value = '121.63'
$ws.Range('E47').Value = '  +2.22%  '
$ws.Range('B48').Value = 'Stellar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.110'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +4.12%  '
$ws.Range('B49').Value = 'Fetch.AI'
$ws.Range('C49').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.00'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +3.36%  '
$ws.Range('D50').Value = '23.62'
$ws.Range('E50').Value = '  +2.72%  '
$ws.Range('D51').Value = '2.02'
$ws.Range('E51').Value = '  +3.22%  '
